# Updated cryptos list on Tue Feb  6 16:42:09 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.322.52"
$ws.Range("E2").Value = "  +1.72%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.354.33"
$ws.Range("E3").Value = "  +2.61%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'303.45"
$ws.Range("E5").Value = "  +0.88%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'95.77"
$ws.Range("E6").Value = "  -0.10%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.504"
$ws.Range("E7").Value = "  -0.23%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.02%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.93%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'34.32"
$ws.Range("E10").Value = "  -0.93%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  +0.32%  "

# Rows 12 & 13 - Chainlink/TRON swap places
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.122"
$ws.Range("E12").Value = "  +2.66%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'18.53"
$ws.Range("E13").Value = "  -3.29%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.11%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.721.10"
$ws.Range("E15").Value = "  +2.61%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.335.75"
$ws.Range("E16").Value = "  +1.74%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  +2.25%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "43.288.57"
$ws.Range("E18").Value = "  +1.73%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("E19").Value = "  +0.74%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  +3.77%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  +0.46%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "'68.32"
$ws.Range("E22").Value = "  +1.17%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "'236.45"
$ws.Range("E23").Value = "  +0.67%  "

# Row 24 - ImmutableX
$ws.Range("D24").Value = "'2.21"
$ws.Range("E24").Value = "  -0.73%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.06%  "

# Row 26 - PancakeSwap
$ws.Range("E26").Value = "  +0.91%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'24.71"
$ws.Range("E27").Value = "  +0.35%  "

# Row 28 - Toncoin
$ws.Range("E28").Value = "  +14.92%  "

# Row 29 - Cosmos
$ws.Range("D29").Value = "'9.18"
$ws.Range("E29").Value = "  +1.79%  "

# Row 30 - InjectiveProtocol
$ws.Range("D30").Value = "'31.56"
$ws.Range("E30").Value = "  -1.94%  "

# Row 31 - FirstDigitalUSD
$ws.Range("E31").Value = "  +0.03%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'5.01"
$ws.Range("E32").Value = "  +1.08%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.0729"
$ws.Range("E33").Value = "  +4.86%  "

# Row 34 - Celestia
$ws.Range("D34").Value = "'17.36"
$ws.Range("E34").Value = "  -0.33%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "'1.84"
$ws.Range("E35").Value = "  +5.49%  "

# Row 36 - RenderToken
$ws.Range("D36").Value = "'4.40"
$ws.Range("E36").Value = "  +0.03%  "

# Row 37 - WEMIXToken
$ws.Range("E37").Value = "  -0.75%  "

# Row 38 - Kaspa
$ws.Range("E38").Value = "  +0.98%  "

# Row 39 - EnergySwap
$ws.Range("D39").Value = "'22.75"
$ws.Range("E39").Value = "  +17.41%  "

# Row 40 - LidoDAOToken
$ws.Range("D40").Value = "'2.76"
$ws.Range("E40").Value = "  +1.31%  "

# Row 41 - Stellar
$ws.Range("E41").Value = "  -0.11%  "

# Row 42 - Monero
$ws.Range("D42").Value = "'114.56"
$ws.Range("E42").Value = "  -30.32%  "

# Row 43 - Maker
$ws.Range("D43").Value = "1.944.63"
$ws.Range("E43").Value = "  -0.54%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  +1.49%  "

# Row 45 - FraxShare
$ws.Range("D45").Value = "'10.02"
$ws.Range("E45").Value = "  -4.48%  "

# Row 46 - ApeXProtocol
$ws.Range("E46").Value = "  +2.41%  "

# Row 47 - NEARProtocol
$ws.Range("E47").Value = "  -0.48%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value = "2.584.72"
$ws.Range("E48").Value = "  +2.55%  "

# Row 49 - MultiversX
$ws.Range("D49").Value = "'53.19"
$ws.Range("E49").Value = "  +0.39%  "

# Row 50 - HuobiToken
$ws.Range("D50").Value = "'2.80"
$ws.Range("E50").Value = "  -3.50%  "

# Row 51 - BitcoinSV
$ws.Range("D51").Value = "'72.29"
$ws.Range("E51").Value = "  +1.01%  "
